# Update the "Max-Min HR Window 10 Sec Stats" sheet with additional test
# cases (rows now run 0..9 through the IF clamp instead of the old
# hand-picked samples) and expose the same stats against a second,
# independent sample set in columns I:K (median over I1:I20 plus a
# dated/blank marker cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Max-Min HR Window 10 Sec Stats")

# --- New sample series in columns A:B (rows 2-12 become 0..10) ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Formula = "=IF(A2<10, A2, 10)"

$ws.Range("A3").Value = 1
$ws.Range("B3").Formula = "=IF(A3<10, A3, 10)"

$ws.Range("A4").Value = 2
$ws.Range("B4").Formula = "=IF(A4<10, A4, 10)"

$ws.Range("A5").Value = 3
$ws.Range("B5").Formula = "=IF(A5<10, A5, 10)"

$ws.Range("A6").Value = 4
$ws.Range("B6").Formula = "=IF(A6<10, A6, 10)"

$ws.Range("A7").Value = 5
$ws.Range("B7").Formula = "=IF(A7<10, A7, 10)"

$ws.Range("A10").Value = 8
$ws.Range("B10").Formula = "=IF(A10<10, A10, 10)"

$ws.Range("A11").Value = 9
$ws.Range("B11").Formula = "=IF(A11<10, A11, 10)"

$ws.Range("A12").Value = 10
$ws.Range("B12").Formula = "=IF(A12<10, A12, 10)"

# --- Stats now computed over the shorter B2:B11 window ---
$ws.Range("C2").Formula = "=MEDIAN(B2:B11)"
$ws.Range("D2").Formula = "=COUNTIF(B2:B11, ""<=""&C2)/COUNT(B2:B11)"
$ws.Range("E2").Formula = "=COUNTIF(B2:B11, "">""&C2)/COUNT(B2:B11)-F2"
$ws.Range("F2").Formula = "=COUNTIF(B2:B11, "">=10"")/COUNT(B2:B11)"

# --- New second sample set exposed in columns I:K ---
$ws.Range("I1").Value = 0
$ws.Range("I2").Value = 10
$ws.Range("J2").Formula = "=MEDIAN(I1:I20)"
$ws.Range("I3").Value = 10

# K2 gets a new number format (date, "d-mmm") but stays otherwise blank -
# this is what introduces the new cellXfs entry in styles.xml.
$ws.Range("K2").NumberFormat = "d-mmm"

# --- Selection moves to the new data entry point ---
$ws.Activate()
$ws.Range("I6").Select()
